# campaign pom class is created
#
# Adds a new "CampaignModule" worksheet (between "product" and the hidden
# "Sheet1") containing the campaign test-data table, and updates the
# selections/active-sheet bookkeeping left behind by the author's editing
# session.

$wb = $excel.ActiveWorkbook

# --- org sheet: selection moves, it is no longer the active tab ---------
$orgSheet = $wb.Worksheets.Item("org")
[void]$orgSheet.Range("C11").Select()

# --- contact sheet: selection moves from A8 to A9 ------------------------
$contactSheet = $wb.Worksheets.Item("contact")
[void]$contactSheet.Range("A9").Select()

# --- insert the new CampaignModule sheet right after "product" ----------
$afterSheet = $wb.Worksheets.Item("product")
$campaignSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$campaignSheet.Name = "CampaignModule"

# Populate in the same order the strings were originally authored so the
# shared-string table is built up in a matching sequence.
$campaignSheet.Range("A1").Value = "TC_ID"
$campaignSheet.Range("B1").Value = "TESTCASENAME"
$campaignSheet.Range("A2").Value = "TC_01"
$campaignSheet.Range("C1").Value = "CAMPAIGN_NAME"
$campaignSheet.Range("D1").Value = "EVENT_NAME"
$campaignSheet.Range("B2").Value = "CreateContactWithEvents"
$campaignSheet.Range("C2").Value = "GreenLiving"
$campaignSheet.Range("D2").Value = "BreezyBrights"

$campaignSheet.Range("A4").Value = "TC_ID"
$campaignSheet.Range("B4").Value = "TESTCASENAME"
$campaignSheet.Range("C4").Value = "CAMPAIGN_NAME"
$campaignSheet.Range("D4").Value = "SEARCHTEXTFIELD"

$campaignSheet.Range("A5").Value = "Tc_02"
$campaignSheet.Range("B5").Value = "CreateCampaignAndDeletingTheCampaign"
$campaignSheet.Range("C5").Value = "GreenLiving"
$campaignSheet.Range("D5").Value = "GreenLiving"

# Column widths matching the authored layout.
$campaignSheet.Columns.Item(1).ColumnWidth = 13.33
$campaignSheet.Columns.Item(2).ColumnWidth = 39.83
$campaignSheet.Columns.Item(3).ColumnWidth = 30.67
$campaignSheet.Columns.Item(4).ColumnWidth = 34.33

# CampaignModule becomes the active sheet/tab with its own selection.
[void]$campaignSheet.Range("D5").Select()

Write-Host "CampaignModule sheet created"
